$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6282.6665
$ws.Range("I116").Value = 10500
$ws.Range("K116").Value = 10500
$ws.Range("M116").Value = -7058

$ws.Range("H132").Value = 2373.0386
$ws.Range("I132").Value = 2565.1738
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 7695.5214
$ws.Range("L132").Value = 2700
$ws.Range("M132").Value = -5165.5214
$ws.Range("N132").Value = -7760

$ws.Range("H134").Value = 41600
$ws.Range("J134").Value = 41600
$ws.Range("L134").Value = 41600
$ws.Range("N134").Value = -51740

$ws.Range("H138").Value = 11238645
$ws.Range("I138").Value = 1263.6666
$ws.Range("J138").Value = 24394602
$ws.Range("K138").Value = 3790.9998
$ws.Range("L138").Value = 73183806
$ws.Range("M138").Value = 1349.0002
$ws.Range("N138").Value = -73194086

$ws.Range("H141").Value = 1869.8636
$ws.Range("I141").Value = 1863.6666
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 5590.9998
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -410.9997999999996
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1692.5834
$ws.Range("I105").Value = 1816.6666
$ws.Range("K105").Value = 1816.6666
$ws.Range("M105").Value = -69.66660000000002

$ws.Range("H107").Value = 1666.3334
$ws.Range("I107").Value = 1499.5
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1499.5
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 420.5
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1667
$ws.Range("I16").Value = 1841.2222
$ws.Range("J16").Value = 1275
$ws.Range("K16").Value = 1841.2222
$ws.Range("L16").Value = 1275
$ws.Range("M16").Value = -1554.2222
$ws.Range("N16").Value = -1849

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H113").Value = 1667
$ws.Range("I113").Value = 1841.2222
$ws.Range("J113").Value = 1275
$ws.Range("K113").Value = 1841.2222
$ws.Range("L113").Value = 1275
$ws.Range("M113").Value = 328.7778000000001
$ws.Range("N113").Value = -5615

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 232.9
$ws.Range("I2").Value = 439.4
$ws.Range("J2").Value = 26.4
$ws.Range("K2").Value = 2636.4
$ws.Range("L2").Value = 158.4
$ws.Range("M2").Value = -2523.4
$ws.Range("N2").Value = -384.4

$ws.Range("H17").Value = 2000
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 6000
$ws.Range("N17").Value = -6338

$ws.Range("H68").Value = 672993.2
$ws.Range("I68").Value = 1983827.4
$ws.Range("J68").Value = 1590.2683
$ws.Range("K68").Value = 5951482.199999999
$ws.Range("L68").Value = 4770.8049
$ws.Range("M68").Value = -5950671.199999999
$ws.Range("N68").Value = -6392.8049

$ws.Range("H71").Value = 672993.2
$ws.Range("I71").Value = 1983827.4
$ws.Range("J71").Value = 1590.2683
$ws.Range("K71").Value = 17854446.6
$ws.Range("L71").Value = 14312.4147
$ws.Range("M71").Value = -17850390.6
$ws.Range("N71").Value = -22424.4147

$ws.Range("H86").Value = 349
$ws.Range("I86").Value = 365.33334
$ws.Range("J86").Value = 300
$ws.Range("K86").Value = 1096.00002
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = 89.99998000000005
$ws.Range("N86").Value = -3272

$ws.Range("H89").Value = 349
$ws.Range("I89").Value = 365.33334
$ws.Range("J89").Value = 300
$ws.Range("K89").Value = 3288.00006
$ws.Range("L89").Value = 2700
$ws.Range("M89").Value = 2639.99994
$ws.Range("N89").Value = -14556

$ws.Range("H112").Value = 1649.56
$ws.Range("J112").Value = 1691.2609
$ws.Range("L112").Value = 5073.7827
$ws.Range("N112").Value = -7289.7827

$ws.Range("H122").Value = 559.5333000000001
$ws.Range("I122").Value = 349.53845
$ws.Range("J122").Value = 1924.5
$ws.Range("K122").Value = 3145.84605
$ws.Range("L122").Value = 17320.5
$ws.Range("M122").Value = -695.8460500000001
$ws.Range("N122").Value = -22220.5

$ws.Range("H125").Value = 1332.5

$ws.Range("H130").Value = 1369.1666
$ws.Range("I130").Value = 1138.3334
$ws.Range("K130").Value = 3415.0002
$ws.Range("M130").Value = 1604.9998

$ws.Range("H133").Value = 2469.5454
$ws.Range("I133").Value = 2861.6667
$ws.Range("J133").Value = 1999
$ws.Range("K133").Value = 8585.000100000001
$ws.Range("L133").Value = 5997
$ws.Range("M133").Value = -3525.000100000001
$ws.Range("N133").Value = -16117

$ws.Range("H134").Value = 32604.188
$ws.Range("I134").Value = 56390.777
$ws.Range("J134").Value = 2021.4286
$ws.Range("K134").Value = 169172.331
$ws.Range("L134").Value = 6064.2858
$ws.Range("M134").Value = -164102.331
$ws.Range("N134").Value = -16204.2858

$ws.Range("H137").Value = 2513.9062
$ws.Range("I137").Value = 1857.5883
$ws.Range("J137").Value = 3257.7334
$ws.Range("K137").Value = 5572.7649
$ws.Range("L137").Value = 9773.200199999999
$ws.Range("M137").Value = -472.7649000000001
$ws.Range("N137").Value = -19973.2002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3600.6
$ws.Range("I80").Value = 3501
$ws.Range("K80").Value = 3501
$ws.Range("M80").Value = -2503

$ws.Range("H83").Value = 3600.6
$ws.Range("I83").Value = 3501
$ws.Range("K83").Value = 17505
$ws.Range("M83").Value = -12513

$ws.Range("H113").Value = 1326.8889
$ws.Range("I113").Value = 1371.5
$ws.Range("J113").Value = 1237.6666
$ws.Range("K113").Value = 1371.5
$ws.Range("L113").Value = 1237.6666
$ws.Range("M113").Value = 798.5
$ws.Range("N113").Value = -5577.6666

$ws.Range("H122").Value = 63749.117
$ws.Range("I122").Value = 79845.06
$ws.Range("J122").Value = 12242.1
$ws.Range("K122").Value = 239535.18
$ws.Range("L122").Value = 36726.3
$ws.Range("M122").Value = -237085.18
$ws.Range("N122").Value = -41626.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1609.25
$ws.Range("I9").Value = 410.57144
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 410.57144
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = -186.57144
$ws.Range("N9").Value = -10448

$ws.Range("H40").Value = 35361.938
$ws.Range("I40").Value = 52820
$ws.Range("J40").Value = 6265.1665
$ws.Range("K40").Value = 52820
$ws.Range("L40").Value = 6265.1665
$ws.Range("M40").Value = -52684
$ws.Range("N40").Value = -6537.1665

$ws.Range("H61").Value = 2875.1738
$ws.Range("I61").Value = 2685.8235
$ws.Range("J61").Value = 3411.6667
$ws.Range("K61").Value = 2685.8235
$ws.Range("L61").Value = 3411.6667
$ws.Range("M61").Value = -2483.8235
$ws.Range("N61").Value = -3815.6667

$ws.Range("H113").Value = 2875.1738
$ws.Range("I113").Value = 2685.8235
$ws.Range("J113").Value = 3411.6667
$ws.Range("K113").Value = 2685.8235
$ws.Range("L113").Value = 3411.6667
$ws.Range("M113").Value = -515.8235
$ws.Range("N113").Value = -7751.6667

$ws.Range("H122").Value = 55560560
$ws.Range("I122").Value = 55560560
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 166681680
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -166679230
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1030.2632
$ws.Range("I113").Value = 1051.5294
$ws.Range("J113").Value = 849.5
$ws.Range("K113").Value = 3154.5882
$ws.Range("L113").Value = 2548.5
$ws.Range("M113").Value = -984.5881999999997
$ws.Range("N113").Value = -6888.5
